# Refresh the hour-on-forecast report figures (rows 2-20, columns B:J)
# with the latest source values so the chart bound to this range no
# longer shows stale (SettingWithCopy-affected) numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 821
$ws.Cells.Item(2, 4).Value = 822
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 17
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 763.0999999999999
$ws.Cells.Item(2, 10).Value = 7.718516577119661

$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 336
$ws.Cells.Item(3, 4).Value = 336
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 9
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 315
$ws.Cells.Item(3, 10).Value = 6.666666666666665

$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 13
$ws.Cells.Item(4, 4).Value = 13
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 11
$ws.Cells.Item(4, 10).Value = 18.18181818181819

$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 302
$ws.Cells.Item(5, 4).Value = 329
$ws.Cells.Item(5, 5).Value = 20
$ws.Cells.Item(5, 6).Value = 7
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 115
$ws.Cells.Item(5, 10).Value = 186.0869565217391

$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 73
$ws.Cells.Item(6, 4).Value = 73
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 53
$ws.Cells.Item(6, 10).Value = 37.73584905660377

$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 63
$ws.Cells.Item(7, 4).Value = 72
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 95
$ws.Cells.Item(7, 10).Value = -24.21052631578947

$ws.Cells.Item(8, 2).Value = 2
$ws.Cells.Item(8, 3).Value = 300
$ws.Cells.Item(8, 4).Value = 322
$ws.Cells.Item(8, 5).Value = 19
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 3
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 38
$ws.Cells.Item(8, 10).Value = 747.3684210526314

$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 27
$ws.Cells.Item(9, 4).Value = 29
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 82
$ws.Cells.Item(9, 10).Value = -64.63414634146342

$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = 325
$ws.Cells.Item(10, 4).Value = 378
$ws.Cells.Item(10, 5).Value = 47
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = 9
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 527
$ws.Cells.Item(10, 10).Value = -28.27324478178368

$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 273
$ws.Cells.Item(11, 4).Value = 275
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 8
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 199
$ws.Cells.Item(11, 10).Value = 38.19095477386936

$ws.Cells.Item(12, 2).Value = 13
$ws.Cells.Item(12, 3).Value = 518
$ws.Cells.Item(12, 4).Value = 617
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).Value = 6
$ws.Cells.Item(12, 8).Value = 70
$ws.Cells.Item(12, 9).Value = 738.4
$ws.Cells.Item(12, 10).Value = -16.4409534127844

$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 28
$ws.Cells.Item(13, 4).Value = 28
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 473
$ws.Cells.Item(13, 10).Value = -94.08033826638477

$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 516
$ws.Cells.Item(14, 4).Value = 554
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 2
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).Value = 33
$ws.Cells.Item(14, 9).Value = 727
$ws.Cells.Item(14, 10).Value = -23.79642365887208

$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = 173
$ws.Cells.Item(15, 4).Value = 177
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 185
$ws.Cells.Item(15, 10).Value = -4.324324324324325

$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = 78
$ws.Cells.Item(17, 4).Value = 80
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 2
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 95
$ws.Cells.Item(17, 10).Value = -15.78947368421053

$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 4
$ws.Cells.Item(18, 10).Value = -25

$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 21
$ws.Cells.Item(19, 4).Value = 21
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 4
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 11
$ws.Cells.Item(19, 10).Value = 90.90909090909092

$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = 26
$ws.Cells.Item(20, 4).Value = 27
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 74
$ws.Cells.Item(20, 10).Value = -63.51351351351351
